$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Cells in column D hold numeric-looking strings but must stay stored
    # as text (inlineStr / shared string), matching the source data feed.
    # Prefixing with an apostrophe forces Excel to keep the literal text;
    # resetting the style back to Normal afterwards strips the
    # quotePrefix formatting flag that the apostrophe entry adds, so the
    # cell's style stays exactly as it was before the edit.
    $range = $ws.Range($addr)
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

Set-TextValue "D2" "247.72"
Set-TextValue "D3" "22.80"
Set-TextValue "D4" "5.297"
Set-TextValue "D5" "0.05726"
Set-TextValue "D6" "3.425"
Set-TextValue "D7" "6.337"
Set-TextValue "D9" "0.8707"
Set-TextValue "D10" "0.1434"
Set-TextValue "D11" "0.07397"
Set-TextValue "D14" "0.09378"
Set-TextValue "D15" "3.884"
Set-TextValue "D16" "0.001576"
Set-TextValue "D17" "0.04821"
Set-TextValue "D18" "0.0005840"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006168"
Set-TextValue "D20" "0.005172"
Set-TextValue "D21" "0.0009965"
Set-TextValue "D23" "3.740"
Set-TextValue "D24" "2.196"
Set-TextValue "D26" "0.1288"
Set-TextValue "D40" "0.03946"
Set-TextValue "D41" "0.006764"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue "D42" "0.1068"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue "D44" "0.008002"
Set-TextValue "D45" "0.00005605"
Set-TextValue "D47" "0.3600"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
